$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F (shifting old F "District" data to G)
$ws.Columns.Item(6).Insert()

# Header
$ws.Cells.Item(2, 6).Value = "Address"

# Address values derived from the B-column (Name + Address, District)
$ws.Cells.Item(3, 6).Value = "SUSRHSGattipuraMagadi"
$ws.Cells.Item(4, 6).Value = "G J C MalebennurHarihar"
$ws.Cells.Item(5, 6).Value = "Bi Bi Raza Girls High School"
$ws.Cells.Item(7, 6).Value = "Bi Bi Raza Girls High School"
$ws.Cells.Item(8, 6).Value = "G E S Kalaghatagi"
$ws.Cells.Item(9, 6).Value = "S G K H S BalaganurSindagi"
$ws.Cells.Item(10, 6).Value = "Govt. P U College For Boys Guledgudd"
$ws.Cells.Item(11, 6).Value = "Dhareppa Katagighan S B H S MalaghanSindagi"
$ws.Cells.Item(12, 6).Value = "G J C H S Junior College Belur Road"
$ws.Cells.Item(14, 6).Value = "KIRAN KUMAR S R G J C N R Pura"
$ws.Cells.Item(15, 6).Value = "S B S Comp. P U CollegeBidaradahally"
$ws.Cells.Item(16, 6).Value = "Govt. High School for BoysR C Road"
$ws.Cells.Item(17, 6).Value = "B E S P U CollegeKadapattiJamkhandi"
$ws.Cells.Item(18, 6).Value = "Govt. P U College (High School Section) Ilkal"
$ws.Cells.Item(19, 6).Value = "G J C (H S) MuthinakoppaN R Pura"
$ws.Cells.Item(22, 6).Value = "Siraj Ul Uloom Urdu High School Manhalli"
$ws.Cells.Item(23, 6).Value = "G H S Pillangere"
$ws.Cells.Item(24, 6).Value = "G H S SangameshwarKalaghatagi"
$ws.Cells.Item(25, 6).Value = "Govt. Girls High School Sedam"
$ws.Cells.Item(26, 6).Value = "S J R Junior College Balehonnur"
$ws.Cells.Item(27, 6).Value = "Poojya Shantaveereshwara High School Sedam"
$ws.Cells.Item(29, 6).Value = "S P P U College (H S) TerdalJamakhandi"
$ws.Cells.Item(31, 6).Value = "G H S Umblebylu"
$ws.Cells.Item(32, 6).Value = "Mudigere"
$ws.Cells.Item(34, 6).Value = "Shri Sharanabasaveshwara High School IddalagiHunagund"
$ws.Cells.Item(35, 6).Value = "S K A R H S Bannikuppe"
$ws.Cells.Item(36, 6).Value = "Sri Guru Karibasaveswara High School UkkadagatriHarihara"
$ws.Cells.Item(38, 6).Value = "K P S Anavatti (G P U C) Soraba"
$ws.Cells.Item(39, 6).Value = "B S B S R H S Bidadi"
$ws.Cells.Item(40, 6).Value = "Sri Lakshman Rao GurjerHigh School JavaliMudigere"
$ws.Cells.Item(41, 6).Value = "S M H S KudurMagadi"
$ws.Cells.Item(42, 6).Value = "Govt. High School Hirebyle Mudigere"
$ws.Cells.Item(43, 6).Value = "G H S MaganageraJewargi"
$ws.Cells.Item(44, 6).Value = "G H S SaluruShikaripura"
$ws.Cells.Item(45, 6).Value = "S V H S Hombegowdana Doddi"
$ws.Cells.Item(46, 6).Value = "S V R High School MotaganahallyMagadi"
$ws.Cells.Item(47, 6).Value = "Banakal High School BanakalMudigere"
$ws.Cells.Item(49, 6).Value = "GuniKalaburagi North"
$ws.Cells.Item(51, 6).Value = "G J C BannikoduHarihara"
$ws.Cells.Item(52, 6).Value = "K P S Anavatti(G P U C) Soraba"
$ws.Cells.Item(53, 6).Value = "S S H S Haramaghatta"
$ws.Cells.Item(54, 6).Value = "Govt. High School Jayawadagim Basavana Bagewadi"
$ws.Cells.Item(55, 6).Value = "Govt. P U College For Girls"
$ws.Cells.Item(56, 6).Value = "Govt. High School C C I KurkuntaSedam"
$ws.Cells.Item(57, 6).Value = "Govt. High School RangolSedam"
